# Week 13 logging: add new RB "De.Jackson" row to the RB sheet,
# and move the active tab/selection from K to RB.

$wb = $excel.ActiveWorkbook
$rb = $wb.Worksheets.Item("RB")

# New player row (row 6) with zeroed weekly stats.
$rb.Range("A6").Value = "De.Jackson"
$rb.Range("B6:J6").Value = 0

# Make RB the active sheet/tab, with the new selection left on J7
# (matches where the cursor lands after entering the last stat).
$rb.Activate() | Out-Null
$rb.Range("J7").Select() | Out-Null
